$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Elimina EC anteriores y se agregan nuevos" - update the Periodo Mora
# column (E16:E20) to the new set of periods, newest first, and move the
# corresponding Valor Mora (F16:F20) figures along with their period.
$ws.Range("E16").Value = "1904"
$ws.Range("E17").Value = "1903"
$ws.Range("E18").Value = "1902"
$ws.Range("E19").Value = "1901"
$ws.Range("E20").Value = "1812"

$ws.Range("F16").Value = 22916
$ws.Range("F17").Value = 31249
$ws.Range("F18").Value = 31249
$ws.Range("F19").Value = 31249
$ws.Range("F20").Value = 31249
